# Adds a new "Profile" worksheet (with Profile test-case data) as the
# first sheet in the workbook, wires up the "Category" defined name that
# points at an external workbook's Data sheet, and refreshes the
# workbookView geometry — matching the authored commit:
# "Added Test cases for Profile features including FullName, Availability,
#  Hours, Earn Target, and Description."

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "Profile" sheet in front of the existing first sheet
# ---------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($firstSheet)
$ws.Name = "Profile"

# ---------------------------------------------------------------------
# 2. Long / special-case text blocks used further down
# ---------------------------------------------------------------------
$longDescription = "This is another test to check if description will accept more than 600 characters. This is another test to check if description will accept more than 600 characters. This is another test to check if description will accept more than 600 characters. This is another test to check if description will accept more than 600 characters. This is another test to check if description will accept more than 600 characters. This is another test to check if description will accept more than 600 characters. This is another test to check if description will accept more than 600 characters. This is another test to check if description will accept more than 600 characters.`nThis is another test to check if description will accept more than 600 characters. This is another test to check if description will accept more than 600 characters. This is another test to check if description will accept more than 600 characters. This is another test to check if description will accept more than 600 characters. This is another test to check if description will accept more than 600 characters. This is another test to check if description will accept more than 600 characters. This is another test to check if description will accept more than 600 characters. This is another test to check if description will accept more than 600 characters.`n"

$leadingSpaceDescription = " This description starts with a space."

$trailingSpacesBase = "Heaps of unwanted spaces after a valid text. Total of 236 characters."
$trailingSpacesPadded = $trailingSpacesBase.PadRight(236 - 1, ' ')
$trailingSpacesDescription = $trailingSpacesPadded + "`n"

$specialCharacters = "!@#$%^&*()_+<>,.?~``"

# ---------------------------------------------------------------------
# 3. Header row
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "FirstName"
$ws.Range("B1").Value = "LastName"
$ws.Range("C1").Value = "FullName"
$ws.Range("D1").Value = "AvailableTime"
$ws.Range("E1").Value = "AvailableHours"
$ws.Range("F1").Value = "EarnTarget"
$ws.Range("G1").Value = "Description"
$ws.Range("H1").Value = "Country"
$ws.Range("I1").Value = "City"
$ws.Range("J1").Value = "Language"
$ws.Range("K1").Value = "Skill"
$ws.Range("L1").Value = "University"
$ws.Range("M1").Value = "Degree"
$ws.Range("N1").Value = "Certificate"
$ws.Range("O1").Value = "CertifiedFrom"
$ws.Range("P1").Value = "Description"

# ---------------------------------------------------------------------
# 4. Row 2 - happy-path profile
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Sheila"
$ws.Range("B2").Value = "Dimasuhid"
$ws.Range("C2").Value = "Sheila Dimasuhid"
$ws.Range("D2").Value = "Full Time"
$ws.Range("E2").Value = "Less than 30hours a week"
$ws.Range("F2").Value = "Less than `$500 per month"
$ws.Range("G2").Value = "I love coding and working on my skill to improve and get better to become a Software Tester."
$ws.Range("H2").Value = "New Zealand"
$ws.Range("I2").Value = "Auckland"
$ws.Range("J2").Value = "English"
$ws.Range("K2").Value = "Automation Testing"
$ws.Range("L2").Value = "Vignan"
$ws.Range("M2").Value = "M.Tech"
$ws.Range("N2").Value = "ISTQB"
$ws.Range("O2").Value = "ANZTB"
$ws.Range("P2").Value = "4 years as Manual Tester and 3 Years as Automation Tester"

# ---------------------------------------------------------------------
# 5. Row 3 - alternate availability / earn-target / description cases
# ---------------------------------------------------------------------
$ws.Range("D3").Value = "Part Time"
$ws.Range("E3").Value = "More than 30hours a week"
$ws.Range("F3").Value = "Between `$500 and `$1000 per month"
$ws.Range("G3").Value = $longDescription
$ws.Range("H3").Value = "jkl"
$ws.Range("I3").Value = "jkl"
$ws.Range("J3").Value = "kl"

# ---------------------------------------------------------------------
# 6. Row 4 - more availability / earn-target / description cases
# ---------------------------------------------------------------------
$ws.Range("E4").Value = "As needed"
$ws.Range("F4").Value = "More than `$1000 per month"
$ws.Range("G4").Value = $specialCharacters

# ---------------------------------------------------------------------
# 7. Row 5 - description leading-space edge case
# ---------------------------------------------------------------------
$ws.Range("G5").Value = $leadingSpaceDescription

# ---------------------------------------------------------------------
# 8. Row 6 - description trailing-spaces edge case
# ---------------------------------------------------------------------
$ws.Range("G6").Value = $trailingSpacesDescription

# ---------------------------------------------------------------------
# 9. Cosmetics: header styling (bold-ish "Normal 2" look carried over
#    from the source workbook the Profile sheet was authored in) and
#    column sizing so the new sheet reads like the rest of the workbook.
# ---------------------------------------------------------------------
$ws.Range("A1:P1").Font.Name = "Calibri"
$ws.Range("A1:P1").Interior.ThemeColor = 5
$ws.Columns.Item(1).ColumnWidth = 9.27
$ws.Columns.Item(2).ColumnWidth = 19.27
$ws.Columns.Item(3).ColumnWidth = 15.18
$ws.Columns.Item(4).ColumnWidth = 20.09
$ws.Columns.Item(5).ColumnWidth = 23.54
$ws.Columns.Item(6).ColumnWidth = 31.45
$ws.Columns.Item(7).ColumnWidth = 33.82
$ws.Columns.Item(8).ColumnWidth = 20.09
$ws.Columns.Item(9).ColumnWidth = 20.09
$ws.Columns.Item(10).ColumnWidth = 18.27
$ws.Columns.Item(11).ColumnWidth = 24.82
$ws.Columns.Item(12).ColumnWidth = 18
$ws.Columns.Item(13).ColumnWidth = 12.73
$ws.Columns.Item(14).ColumnWidth = 14.45
$ws.Columns.Item(15).ColumnWidth = 18.45
$ws.Columns.Item(16).ColumnWidth = 52.73

$ws.Range("A2").Select()
$ws.Application.ActiveWindow.ScrollRow = 2

# ---------------------------------------------------------------------
# 10. Defined name "Category" pointing at the external "Data" sheet.
# ---------------------------------------------------------------------
$wb.Names.Add("Category", "=[1]Data!`$B`$4:`$B`$11")

$wb.Worksheets.Item("Profile").Select()
